$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in main block assignments: the "C# B#- ..." condition strings
# in column B (rows 2-7) were leftover/incorrect and should be removed,
# leaving column B populated only in the header row.
$ws.Range("B2:B7").ClearContents()

# Update the active selection and zoom level to reflect the reviewed state
[void]$ws.Range("B1").Select()
$ws.Application.ActiveWindow.Zoom = 125
